# New crime data collected - weekly CompStat update
# Updates the report header (volume number, week-covering dates) and the
# Crime Complaints table (rows 14-27, columns C:N) for the 1st Precinct
# weekly report.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Header text: "Volume 31   Number  6" -> "...7" and the week-covering
#    date range "2/5/2024 ... 2/11/2024" -> "2/12/2024 ... 2/18/2024"
# ---------------------------------------------------------------------
$ws.Range("A8").Value = "Volume 31   Number  7"
$ws.Range("C9").Value = "Report Covering the Week  2/12/2024  Through  2/18/2024"

# ---------------------------------------------------------------------
# 2. Cells that flip between a numeric value and the text placeholder
#    "0" (shared string used whenever a count is zero). Copying the
#    already-correctly-styled neighbour cell keeps font/alignment/number
#    format identical to the rest of the table, then the value is set
#    explicitly on top of the copy.
# ---------------------------------------------------------------------

# Row 14 (Murder): Week-to-date 2024 count goes from 1 to 0
$ws.Range("D14").Copy($ws.Range("C14"))
$ws.Range("C14").Value = "0"

# Row 15 (Rape): Week-to-date 2024 count goes from 2 to 0
$ws.Range("D15").Copy($ws.Range("C15"))
$ws.Range("C15").Value = "0"

# Row 20 (G.L.A.): Week-to-date 2024 goes from 1 to "0"; 2023 goes from
# "0" to a real count of 2, which also makes the %Chg (E20) and 2-Year
# %Chg (M20) columns become real numbers instead of the "***.*" text.
$ws.Range("D20").Copy($ws.Range("C20"))
$ws.Range("C20").Value = "0"
$ws.Range("G20").Copy($ws.Range("D20"))
$ws.Range("D20").Value = 2
$ws.Range("H20").Copy($ws.Range("E20"))
$ws.Range("E20").Value = -100
$ws.Range("H20").Copy($ws.Range("M20"))
$ws.Range("M20").Value = 50

# Row 22 (Transit): Week-to-date 2024 goes from "0" to a real count of 2
$ws.Range("D22").Copy($ws.Range("C22"))
$ws.Range("C22").Value = 2

# Row 26 (UCR Rape*): Week-to-date 2024 goes from 2 to "0"
$ws.Range("D26").Copy($ws.Range("C26"))
$ws.Range("C26").Value = "0"

# ---------------------------------------------------------------------
# 3. Remaining numeric-only cell updates (no type/style change) across
#    the Crime Complaints table.
# ---------------------------------------------------------------------
$updates = @(
    @("C16", 1),
    @("E16", -66.666666666666),
    @("G16", 13),
    @("H16", -23.076923076923),
    @("J16", 21),
    @("K16", -23.809523809523),
    @("L16", -27.272727272727),
    @("N16", -87.878787878787),
    @("D17", 2),
    @("E17", 0),
    @("F17", 8),
    @("H17", 0),
    @("I17", 13),
    @("J17", 15),
    @("K17", -13.333333333333),
    @("L17", 8.333333333333),
    @("M17", 44.444444444444),
    @("N17", -31.578947368421),
    @("C18", 11),
    @("D18", 3),
    @("E18", 266.666666666667),
    @("F18", 19),
    @("G18", 18),
    @("H18", 5.555555555555),
    @("I18", 30),
    @("J18", 25),
    @("K18", 20),
    @("L18", -9.090909090909),
    @("M18", 0),
    @("N18", -73.913043478260),
    @("C19", 24),
    @("D19", 25),
    @("E19", -4),
    @("F19", 77),
    @("H19", -6.097560975609),
    @("I19", 132),
    @("J19", 147),
    @("K19", -10.204081632653),
    @("L19", -20),
    @("M19", -6.382978723404),
    @("N19", -71.052631578947),
    @("G20", 3),
    @("H20", -66.666666666666),
    @("J20", 6),
    @("K20", -50),
    @("L20", 50),
    @("N20", -97.674418604651),
    @("C21", 38),
    @("D21", 35),
    @("E21", 8.571428571428),
    @("F21", 118),
    @("G21", 126),
    @("H21", -6.349206349206),
    @("I21", 198),
    @("J21", 216),
    @("K21", -8.333333333333),
    @("L21", -16.101694915254),
    @("M21", 4.210526315789),
    @("N21", -76.760563380281),
    @("D22", 3),
    @("E22", -33.333333333333),
    @("F22", 5),
    @("G22", 10),
    @("H22", -50),
    @("I22", 8),
    @("J22", 16),
    @("K22", -50),
    @("L22", -38.461538461538),
    @("M22", -38.461538461538),
    @("C24", 67),
    @("D24", 63),
    @("E24", 6.349206349206),
    @("F24", 278),
    @("G24", 282),
    @("H24", -1.418439716312),
    @("I24", 518),
    @("J24", 497),
    @("K24", 4.225352112676),
    @("L24", 7.024793388429),
    @("M24", 168.39378238342),
    @("C25", 10),
    @("D25", 7),
    @("E25", 42.857142857142),
    @("F25", 35),
    @("G25", 23),
    @("H25", 52.173913043478),
    @("I25", 49),
    @("J25", 40),
    @("K25", 22.5),
    @("L25", 25.641025641025),
    @("M25", 63.333333333333),
    @("D27", 1),
    @("E27", 0),
    @("G27", 9),
    @("H27", -11.111111111111),
    @("I27", 11),
    @("J27", 13),
    @("K27", -15.384615384615),
    @("L27", 37.5)
)

foreach ($u in $updates) {
    $ws.Range($u[0]).Value = $u[1]
}

Write-Host "Applied weekly crime data update."
